# Logboek.xlsx edit: add "Strategie 5" entries (two new log rows: 15 and 16)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: "Latex document aanpassen" on 2015-04-11 (x for Jari), time 30min ---
$ws.Range("A15").Value = "Latex document aanpassen"
$ws.Range("E15").Value = "x"
$ws.Range("G15").Value = "30min"

# Date cell for row 15 gets a new, plain (non-bold-look) date style.
$ws.Range("B15").NumberFormat = "mm-dd-yy"
$ws.Range("B15").Value = 42105

# --- Row 16: "Simulatie gegevens nakijken + verbeteren " on 2015-04-12 (x for Jari), time 1u ---
$ws.Range("A16").Value = "Simulatie gegevens nakijken + verbeteren "
$ws.Range("E16").Value = "x"
$ws.Range("G16").Value = "1u"

# Date cell for row 16 reuses the bold-date style already used for the other
# log dates (e.g. B14), by copying that cell's format.
$ws.Range("B14").Copy($ws.Range("B16"))
$ws.Range("B16").Value = 42106

# Update the active selection as recorded in the saved sheet view.
$ws.Range("F23").Select()
